$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting existing data rows (2..51) down to (3..52)
$ws.Rows.Item(2).Insert()

# The inserted row inherits formatting from the row above (the header); clear it
# back to the plain (unstyled) look used by the rest of the data rows.
$ws.Rows.Item(2).ClearFormats()

# Populate the new row 2 with the latest week's record
$ws.Range("A2").Value = 5
$ws.Range("B2").Value = "Macroferia Regional de Talca"
$ws.Range("C2").Value = "Maule"
$ws.Range("D2").Value = 44956
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 100112043
$ws.Range("G2").Value = "Pepino dulce"
$ws.Range("H2").Value = "Cultivar IV Región"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 18000
$ws.Range("L2").Value = 18000
$ws.Range("M2").Value = 18000
$ws.Range("N2").Value = "`$/bandeja 18 kilos"
$ws.Range("O2").Value = "Provincia de Limarí"
$ws.Range("P2").Value = 1000
$ws.Range("Q2").Value = 18
$ws.Range("R2").Value = "Hortaliza"

# D column holds dates stored as serials with a custom date number format;
# re-apply it to the new row's D cell (same format used by the rest of column D).
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
